$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46032
$ws.Range("B2").Value = 75.17
$ws.Range("C2").Value = 70.12
$ws.Range("D2").Value = 65.22
$ws.Range("E2").Value = 59.91
$ws.Range("F2").Value = 55
$ws.Range("G2").Value = 54.27
$ws.Range("H2").Value = 62.89
$ws.Range("I2").Value = 70.3
$ws.Range("J2").Value = 80.47
$ws.Range("K2").Value = 83.31999999999999
$ws.Range("L2").Value = 64.98999999999999
$ws.Range("M2").Value = 47.14
$ws.Range("N2").Value = 28.77
$ws.Range("O2").Value = 30.79
$ws.Range("P2").Value = 36.83
$ws.Range("Q2").Value = 48.21
$ws.Range("R2").Value = 71.22
$ws.Range("S2").Value = 96.14
$ws.Range("T2").Value = 118.54
$ws.Range("U2").Value = 126.26
$ws.Range("V2").Value = 112.26
$ws.Range("W2").Value = 101.17
$ws.Range("X2").Value = 100.43
$ws.Range("Y2").Value = 98.2
$ws.Range("Z2").Value = 73.23
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 103.04
$ws.Range("AD2").Value = 122.4
$ws.Range("AF2").Value = 106.72
